$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '37.218.91'
$c.Style = 'Normal'
$ws.Range("E2").Value = '  +0.16%  '
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '2.056.24'
$c.Style = 'Normal'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  -0.14%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '248.84'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  -1.95%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '0.665'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  -1.94%  '
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '57.93'
$c.Style = 'Normal'
$ws.Range("E7").Value = '  -1.55%  '
$ws.Range("E8").Value = '  -0.02%  '
$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.383'
$c.Style = 'Normal'
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.0780'
$c.Style = 'Normal'
$ws.Range("E10").Value = '  -2.63%  '
$ws.Range("E11").Value = '  +0.21%  '
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '16.07'
$c.Style = 'Normal'
$ws.Range("E12").Value = '  -1.45%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.883'
$c.Style = 'Normal'
$ws.Range("E13").Value = '  +7.26%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '2.353.93'
$c.Style = 'Normal'
$ws.Range("E14").Value = '  -1.00%  '
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '5.72'
$c.Style = 'Normal'
$ws.Range("E15").Value = '  +3.30%  '
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '2.053.47'
$c.Style = 'Normal'
$ws.Range("E16").Value = '  -1.06%  '
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '18.32'
$c.Style = 'Normal'
$ws.Range("E17").Value = '  +16.11%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '37.219.52'
$c.Style = 'Normal'
$ws.Range("E18").Value = '  +0.14%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '75.00'
$c.Style = 'Normal'
$ws.Range("E19").Value = '  +0.33%  '
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '0.0₃0896'
$c.Style = 'Normal'
$ws.Range("E20").Value = '  -3.36%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '5.41'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  -1.60%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '237.60'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("E23").Value = '  +0.05%  '
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '2.48'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  +3.31%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '2.20'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  -3.84%  '
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '9.52'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  +1.89%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '169.73'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  -0.05%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '20.20'
$c.Style = 'Normal'
$ws.Range("E28").Value = '  -0.93%  '
$ws.Range("E29").Value = '  -1.25%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '4.85'
$c.Style = 'Normal'
$ws.Range("E30").Value = '  +0.91%  '
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '1.14'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  -1.21%  '
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '0.0621'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  -2.13%  '
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '4.49'
$c.Style = 'Normal'
$ws.Range("E33").Value = '  -0.70%  '
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '0.0893'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  -1.96%  '
$ws.Range("E35").Value = '  +0.00%  '
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.Style = 'Normal'
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("E38").Value = '  -1.73%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '5.27'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  +15.10%  '
$ws.Range("E40").Value = '  +11.03%  '
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '0.0995'
$c.Style = 'Normal'
$ws.Range("E41").Value = '  -14.72%  '
$ws.Range("E42").Value = '  -1.78%  '
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '17.36'
$c.Style = 'Normal'
$ws.Range("E43").Value = '  -2.64%  '
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '1.15'
$c.Style = 'Normal'
$ws.Range("E44").Value = '  -2.54%  '
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '96.39'
$c.Style = 'Normal'
$ws.Range("E45").Value = '  -2.81%  '
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.Style = 'Normal'
$ws.Range("E46").Value = '  -1.52%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '1.273.27'
$c.Style = 'Normal'
$ws.Range("E47").Value = '  -2.61%  '
$ws.Range("E48").Value = '  -2.76%  '
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '6.83'
$c.Style = 'Normal'
$ws.Range("E49").Value = '  -1.53%  '
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '2.242.44'
$c.Style = 'Normal'
$ws.Range("E50").Value = '  -0.85%  '
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '44.04'
$c.Style = 'Normal'
$ws.Range("E51").Value = '  -0.80%  '